# Update cryptos list with latest scraped price/volume data
# (commit: "Updated cryptos list on Thu Nov  7 22:39:21 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. All of these columns (Coin/Link/Price/
# Volume) are plain text cells, so NumberFormat is forced to "@" (Text) while the
# value is written to stop Excel from re-interpreting things like "195.20" or
# "2.873.74" as numbers/dates, then the style is reset back to Normal (matching
# the original, unstyled data cells) once the literal text is safely stored.
$updates = [ordered]@{
    'D2' = '75.729.62'
    'E2' = '  +0.12%  '
    'D3' = '2.873.74'
    'E3' = '  +6.46%  '
    'E4' = '  +0.00%  '
    'D5' = '195.20'
    'E5' = '  +3.14%  '
    'D6' = '597.30'
    'E6' = '  +1.34%  '
    'E7' = '  +0.05%  '
    'D8' = '0.553'
    'E8' = '  +2.13%  '
    'D9' = '0.192'
    'E9' = '  -2.64%  '
    'D10' = '2.878.79'
    'E10' = '  +6.56%  '
    'D11' = '0.400'
    'E11' = '  +10.96%  '
    'E12' = '  -1.53%  '
    'D13' = '4.91'
    'E13' = '  +3.34%  '
    'D14' = '3.426.02'
    'E14' = '  +7.63%  '
    'D15' = '75.727.56'
    'E15' = '  +0.35%  '
    'D16' = '0.0000189'
    'E16' = '  -0.15%  '
    'D17' = '27.31'
    'E17' = '  +2.19%  '
    'D18' = '2.899.03'
    'E18' = '  +7.94%  '
    'D19' = '8.94'
    'E19' = '  -5.56%  '
    'D20' = '12.56'
    'E20' = '  +3.37%  '
    'D21' = '381.28'
    'E21' = '  +0.60%  '
    'D22' = '2.30'
    'E22' = '  -0.21%  '
    'D23' = '4.14'
    'E23' = '  +1.48%  '
    'D24' = '71.65'
    'E24' = '  +1.26%  '
    'E25' = '  -0.01%  '
    'D27' = '4.23'
    'E27' = '  +0.41%  '
    'D28' = '9.74'
    'E28' = '  +2.12%  '
    'D29' = '0.0000107'
    'E29' = '  +10.71%  '
    'E30' = '  -0.25%  '
    'D31' = '1.40'
    'E31' = '  -1.31%  '
    'D32' = '507.72'
    'E32' = '  -2.83%  '
    'D33' = '7.77'
    'E33' = '  -1.20%  '
    'D34' = '1.82'
    'E34' = '  +2.26%  '
    'E35' = '  +0.05%  '
    'D36' = '164.67'
    'E36' = '  +1.33%  '
    'D37' = '20.13'
    'E37' = '  +3.68%  '
    'D38' = '19.66'
    'E38' = '  +1.46%  '
    'E39' = '  -5.45%  '
    'D40' = '183.20'
    'E40' = '  +5.73%  '
    'E41' = '  -0.11%  '
    'D42' = '0.344'
    'E42' = '  +3.36%  '
    'D43' = '5.00'
    'E43' = '  -1.19%  '
    'D44' = '1.68'
    'E44' = '  -2.16%  '
    'D45' = '0.0909'
    'E45' = '  +6.73%  '
    'D46' = '1.22'
    'E46' = '  +1.07%  '
    'D47' = '40.22'
    'E47' = '  +2.20%  '
    'D48' = '2.36'
    'E48' = '  -2.10%  '
    'E49' = '  +6.13%  '
    'B50' = 'Mantle'
    'C50' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D50' = '0.668'
    'E50' = '  +12.18%  '
    'B51' = 'Filecoin'
    'C51' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D51' = '3.77'
    'E51' = '  +2.20%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
